$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 12 new quarterly rows (2010-Q1 .. 2012-Q4) above the existing
# 2013-Q1 row, pushing the rest of the table down by 12 rows.
$ws.Range("A2:C13").EntireRow.Insert()

# New data: Time Period, Gross domestic product (AUD, Millions), GDP/Capita
$data = @(
  @("2010-Q1", 314838, 14336),
  @("2010-Q2", 340575, 15460),
  @("2010-Q3", 345512, 15633),
  @("2010-Q4", 365403, 16482),
  @("2011-Q1", 341094, 15319),
  @("2011-Q2", 366641, 16414),
  @("2011-Q3", 372918, 16626),
  @("2011-Q4", 388419, 17248),
  @("2012-Q1", 358492, 15836),
  @("2012-Q2", 381023, 16763),
  @("2012-Q3", 379954, 16642),
  @("2012-Q4", 396614, 17301)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = 2 + $i
  $ws.Cells.Item($row, 1).Value = $data[$i][0]
  $ws.Cells.Item($row, 2).Value = $data[$i][1]
  $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Match the centered formatting used for the new rows.
$ws.Range("A2:C13").Font.Bold = $false
$ws.Range("A2:C13").HorizontalAlignment = -4108

# Restore the selected cell to match the saved view state.
$ws.Range("J13").Select()
